$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "number of units" (G column) inputs that drive the billing
# formulas in column I. These were previously blank (0 units -> 0 amount).
$ws.Range("G9").Value  = 116
$ws.Range("G12").Value = 116
$ws.Range("G14").Value = 118
$ws.Range("G16").Value = 27
$ws.Range("G17").Value = 118
$ws.Range("G18").Value = 118
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1

$wb.Application.CalculateFull()
